$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.801.28'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +1.51%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.149.43'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +0.15%  '

$ws.Range("E4").Value = '  +0.46%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '592.33'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.10%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '153.32'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +4.82%  '

$ws.Range("E7").Value = '  +0.15%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.146.20'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +0.45%  '

$ws.Range("E9").Value = '  +1.13%  '

$ws.Range("E10").Value = '  +0.91%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.02'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +4.65%  '

$ws.Range("E12").Value = '  +1.85%  '

$ws.Range("B13").Value = 'ShibaInu'
$ws.Range("C13").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000248'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +0.54%  '

$ws.Range("B14").Value = 'Avalanche'
$ws.Range("C14").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '38.46'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +4.12%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.673.90'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +0.25%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.119'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -1.49%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.30'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +3.24%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '64.373.94'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +1.13%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.149.46'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +0.18%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '476.47'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +2.79%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '15.04'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +5.38%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.748'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +2.09%  '

$ws.Range("E23").Value = '  +3.63%  '

$ws.Range("B24").Value = 'Fetch.AI'
$ws.Range("C24").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.43'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +10.29%  '

$ws.Range("B25").Value = 'InternetComputer(DFINITY)'
$ws.Range("C25").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '13.45'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +4.24%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '82.44'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +1.69%  '

$ws.Range("B27").Value = 'RenderToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.01'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +7.21%  '

$ws.Range("B28").Value = 'Dai'
$ws.Range("C28").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.09%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.46'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +4.91%  '

$ws.Range("E30").Value = '  +1.97%  '

$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.24'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +1.42%  '

$ws.Range("B32").Value = 'FirstDigitalUSD'
$ws.Range("C32").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.01'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +0.56%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.119'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +9.03%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '27.77'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +3.33%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0₃0879'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +3.64%  '

$ws.Range("B36").Value = 'Mantle'
$ws.Range("C36").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.07'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +2.93%  '

$ws.Range("B37").Value = 'dogwifhat'
$ws.Range("C37").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.45'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +3.24%  '

$ws.Range("E38").Value = '  +3.79%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.32'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +0.50%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '469.27'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +7.62%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '9.43'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +6.44%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '51.32'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.26%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.299'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +6.96%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0380'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +2.49%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.897.95'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.55%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.111'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +3.18%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '39.53'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +6.43%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '131.83'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +4.13%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '25.91'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +7.58%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.32'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +6.66%  '

$ws.Range("E51").Value = '  +0.04%  '
